$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(60, 1).Value = "I would like to have my hair cut."
$ws.Cells.Item(60, 2).Value = "カットをお願いします。|カットをおねがいします。"
$ws.Cells.Item(61, 1).Value = "Please don't make it too short."
$ws.Cells.Item(61, 2).Value = "あまり短くしないでください。|あまりみじかくしないでください。"
$ws.Cells.Item(62, 1).Value = "Please don't shave me."
$ws.Cells.Item(62, 2).Value = "そらないでください。"
$ws.Cells.Item(63, 1).Value = "Please cut off about 3 centimeters."
$ws.Cells.Item(63, 2).Value = "３センチぐらい切ってください。|３センチぐらいきってください。"
$ws.Cells.Item(64, 1).Value = "Please cut the back all the same length."
$ws.Cells.Item(64, 2).Value = "後ろをそろえてください。|うしろをそろえてください。"
$ws.Cells.Item(65, 1).Value = "Please dye my hair red."
$ws.Cells.Item(65, 2).Value = "赤にそめてください。|あかにそめてください。"
$ws.Cells.Item(66, 1).Value = "I want my hair to be like Bob Marley's."
$ws.Cells.Item(66, 2).Value = "ボブ・マーリーみたいな髪形にしたいんですが。|ボブ・マーリーみたいなかみがたにしたいんですが。"
$ws.Cells.Item(67, 1).Value = "shampoo"
$ws.Cells.Item(67, 2).Value = "シャンプー"
$ws.Cells.Item(68, 1).Value = "cut"
$ws.Cells.Item(68, 2).Value = "カット"
$ws.Cells.Item(69, 1).Value = "blow-dry"
$ws.Cells.Item(69, 2).Value = "ブロー"
$ws.Cells.Item(70, 1).Value = "perm"
$ws.Cells.Item(70, 2).Value = "パーマ"
$ws.Cells.Item(71, 1).Value = "hair coloring"
$ws.Cells.Item(71, 2).Value = "カラー"
$ws.Cells.Item(72, 1).Value = "set"
$ws.Cells.Item(72, 2).Value = "セット"
$ws.Cells.Item(73, 1).Value = "hair style"
$ws.Cells.Item(73, 2).Value = "髪形|かみがた"
$ws.Cells.Item(74, 1).Value = "to cut"
$ws.Cells.Item(74, 2).Value = "切る|きる"
$ws.Cells.Item(75, 1).Value = "to shave"
$ws.Cells.Item(75, 2).Value = "そる"
$ws.Cells.Item(76, 1).Value = "to crop"
$ws.Cells.Item(76, 2).Value = "刈る|かる"
$ws.Cells.Item(77, 1).Value = "to dye"
$ws.Cells.Item(77, 2).Value = "そめる"
$ws.Cells.Item(78, 1).Value = "to make hair even; to trim"
$ws.Cells.Item(78, 2).Value = "そろえる"
$ws.Cells.Item(79, 1).Value = "to have one's hair permed"
$ws.Cells.Item(79, 2).Value = "パーマをかける"
$ws.Cells.Item(80, 1).Value = "to thin out (hair)"
$ws.Cells.Item(80, 2).Value = "すく"
$ws.Cells.Item(81, 1).Value = "parting (of the hair)"
$ws.Cells.Item(81, 2).Value = "分け目|わけめ"
$ws.Cells.Item(82, 1).Value = "bangs"
$ws.Cells.Item(82, 2).Value = "前髪|まえがみ"
$ws.Cells.Item(83, 1).Value = "side"
$ws.Cells.Item(83, 2).Value = "横|よこ"
$ws.Cells.Item(84, 1).Value = "back"
$ws.Cells.Item(84, 2).Value = "後ろ|うしろ"
